$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-03 12:47:31"

$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
